$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33
$ws.Cells.Item(33, 2).Value = 6845250
$ws.Cells.Item(33, 5).Value = 'Arbroath'
$ws.Cells.Item(33, 6).Value = 'Inverness CT'
$ws.Cells.Item(33, 7).Value = 2
$ws.Cells.Item(33, 8).Value = 3
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 10).Value = 2
$ws.Cells.Item(33, 11).Value = 'A'
$ws.Cells.Item(33, 13).Value = 3.4
$ws.Cells.Item(33, 14).Value = 3.25
$ws.Cells.Item(33, 15).Value = 2.25
$ws.Cells.Item(33, 17).Value = 2.9
$ws.Cells.Item(33, 19).Value = 1.975
$ws.Cells.Item(33, 20).Value = 1.825
$ws.Cells.Item(33, 21).Value = 2.25
$ws.Cells.Item(33, 22).Value = 1.775
$ws.Cells.Item(33, 23).Value = 2.025
$ws.Cells.Item(33, 24).Value = -1
$ws.Cells.Item(33, 26).Value = 1.9
$ws.Cells.Item(33, 27).Value = -1
$ws.Cells.Item(33, 28).Value = 0.825
$ws.Cells.Item(33, 29).Value = 0.7749999999999999
$ws.Cells.Item(33, 30).Value = -1

# Row 34
$ws.Cells.Item(34, 2).Value = 6845251
$ws.Cells.Item(34, 5).Value = 'Ayr'
$ws.Cells.Item(34, 6).Value = 'Raith'
$ws.Cells.Item(34, 7).Value = 1
$ws.Cells.Item(34, 8).Value = 2
$ws.Cells.Item(34, 12).Value = 3.25
$ws.Cells.Item(34, 14).Value = 2
$ws.Cells.Item(34, 15).Value = 2.8
$ws.Cells.Item(34, 16).Value = 3.25
$ws.Cells.Item(34, 17).Value = 2.3
$ws.Cells.Item(34, 18).Value = 0.25
$ws.Cells.Item(34, 19).Value = 1.8
$ws.Cells.Item(34, 20).Value = 2.05
$ws.Cells.Item(34, 21).Value = 2.5
$ws.Cells.Item(34, 22).Value = 2.05
$ws.Cells.Item(34, 23).Value = 1.8
$ws.Cells.Item(34, 26).Value = 1.3
$ws.Cells.Item(34, 28).Value = 1.05
$ws.Cells.Item(34, 29).Value = 1.05

# Row 35
$ws.Cells.Item(35, 2).Value = 6845252
$ws.Cells.Item(35, 5).Value = 'Dundee Utd'
$ws.Cells.Item(35, 6).Value = 'Queens Park'
$ws.Cells.Item(35, 7).Value = 4
$ws.Cells.Item(35, 8).Value = 1
$ws.Cells.Item(35, 9).Value = 1
$ws.Cells.Item(35, 10).Value = 1
$ws.Cells.Item(35, 11).Value = 'H'
$ws.Cells.Item(35, 12).Value = 1.7
$ws.Cells.Item(35, 13).Value = 3.6
$ws.Cells.Item(35, 14).Value = 4.2
$ws.Cells.Item(35, 15).Value = 1.533
$ws.Cells.Item(35, 16).Value = 4
$ws.Cells.Item(35, 17).Value = 5.5
$ws.Cells.Item(35, 18).Value = -1
$ws.Cells.Item(35, 19).Value = 1.875
$ws.Cells.Item(35, 20).Value = 1.975
$ws.Cells.Item(35, 21).Value = 3
$ws.Cells.Item(35, 22).Value = 1.975
$ws.Cells.Item(35, 23).Value = 1.875
$ws.Cells.Item(35, 24).Value = 0.5329999999999999
$ws.Cells.Item(35, 26).Value = -1
$ws.Cells.Item(35, 27).Value = 0.875
$ws.Cells.Item(35, 28).Value = -1
$ws.Cells.Item(35, 29).Value = 0.9750000000000001

# Row 36
$ws.Cells.Item(36, 2).Value = 6845253
$ws.Cells.Item(36, 5).Value = 'Morton'
$ws.Cells.Item(36, 6).Value = 'Airdrieonians'
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 'A'
$ws.Cells.Item(36, 12).Value = 2
$ws.Cells.Item(36, 13).Value = 3.25
$ws.Cells.Item(36, 14).Value = 3.4
$ws.Cells.Item(36, 15).Value = 2.2
$ws.Cells.Item(36, 16).Value = 3.2
$ws.Cells.Item(36, 17).Value = 3.1
$ws.Cells.Item(36, 18).Value = -0.25
$ws.Cells.Item(36, 19).Value = 1.95
$ws.Cells.Item(36, 20).Value = 1.85
$ws.Cells.Item(36, 21).Value = 2.75
$ws.Cells.Item(36, 22).Value = 1.925
$ws.Cells.Item(36, 24).Value = -1
$ws.Cells.Item(36, 26).Value = 2.1
$ws.Cells.Item(36, 27).Value = -1
$ws.Cells.Item(36, 28).Value = 0.8500000000000001
$ws.Cells.Item(36, 29).Value = -1
$ws.Cells.Item(36, 30).Value = 0.875

# Row 37
$ws.Cells.Item(37, 2).Value = 6845254
$ws.Cells.Item(37, 5).Value = 'Partick'
$ws.Cells.Item(37, 6).Value = 'Dunfermline'
$ws.Cells.Item(37, 7).Value = 3
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 1
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 'H'
$ws.Cells.Item(37, 15).Value = 2.3
$ws.Cells.Item(37, 17).Value = 2.875
$ws.Cells.Item(37, 19).Value = 2
$ws.Cells.Item(37, 20).Value = 1.8
$ws.Cells.Item(37, 22).Value = 2
$ws.Cells.Item(37, 23).Value = 1.8
$ws.Cells.Item(37, 24).Value = 1.3
$ws.Cells.Item(37, 26).Value = -1
$ws.Cells.Item(37, 27).Value = 1
$ws.Cells.Item(37, 28).Value = -1
$ws.Cells.Item(37, 29).Value = 0.5
$ws.Cells.Item(37, 30).Value = -0.5

# Row 38
$ws.Cells.Item(38, 2).Value = 6845255
$ws.Cells.Item(38, 5).Value = 'Airdrieonians'
$ws.Cells.Item(38, 6).Value = 'Ayr'
$ws.Cells.Item(38, 7).Value = 1
$ws.Cells.Item(38, 8).Value = 2
$ws.Cells.Item(38, 10).Value = 1
$ws.Cells.Item(38, 11).Value = 'A'
$ws.Cells.Item(38, 12).Value = 2.2
$ws.Cells.Item(38, 13).Value = 3.4
$ws.Cells.Item(38, 14).Value = 2.75
$ws.Cells.Item(38, 15).Value = 2.05
$ws.Cells.Item(38, 16).Value = 3.4
$ws.Cells.Item(38, 17).Value = 2.9
$ws.Cells.Item(38, 18).Value = -0.25
$ws.Cells.Item(38, 19).Value = 1.9
$ws.Cells.Item(38, 20).Value = 1.95
$ws.Cells.Item(38, 21).Value = 2.5
$ws.Cells.Item(38, 22).Value = 1.85
$ws.Cells.Item(38, 23).Value = 2
$ws.Cells.Item(38, 25).Value = -1
$ws.Cells.Item(38, 26).Value = 1.9
$ws.Cells.Item(38, 27).Value = -1
$ws.Cells.Item(38, 28).Value = 0.95
$ws.Cells.Item(38, 29).Value = 0.8500000000000001
$ws.Cells.Item(38, 30).Value = -1

# Row 39
$ws.Cells.Item(39, 2).Value = 6845257
$ws.Cells.Item(39, 5).Value = 'Inverness CT'
$ws.Cells.Item(39, 6).Value = 'Partick'
$ws.Cells.Item(39, 12).Value = 3
$ws.Cells.Item(39, 13).Value = 3.25
$ws.Cells.Item(39, 14).Value = 2.1
$ws.Cells.Item(39, 15).Value = 2.6
$ws.Cells.Item(39, 17).Value = 2.375
$ws.Cells.Item(39, 18).Value = 0
$ws.Cells.Item(39, 22).Value = 1.9
$ws.Cells.Item(39, 23).Value = 1.95
$ws.Cells.Item(39, 27).Value = 0
$ws.Cells.Item(39, 28).Value = 0
$ws.Cells.Item(39, 30).Value = 0.95

# Row 40
$ws.Cells.Item(40, 2).Value = 6845258
$ws.Cells.Item(40, 5).Value = 'Queens Park'
$ws.Cells.Item(40, 6).Value = 'Morton'
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 'D'
$ws.Cells.Item(40, 12).Value = 1.833
$ws.Cells.Item(40, 14).Value = 3.6
$ws.Cells.Item(40, 15).Value = 2.2
$ws.Cells.Item(40, 16).Value = 3.3
$ws.Cells.Item(40, 17).Value = 2.8
$ws.Cells.Item(40, 19).Value = 2.025
$ws.Cells.Item(40, 20).Value = 1.825
$ws.Cells.Item(40, 21).Value = 2.75
$ws.Cells.Item(40, 22).Value = 1.825
$ws.Cells.Item(40, 23).Value = 2.025
$ws.Cells.Item(40, 25).Value = 2.3
$ws.Cells.Item(40, 26).Value = -1
$ws.Cells.Item(40, 27).Value = -0.5
$ws.Cells.Item(40, 28).Value = 0.4125
$ws.Cells.Item(40, 29).Value = -1
$ws.Cells.Item(40, 30).Value = 1.025

# Row 50
$ws.Cells.Item(50, 2).Value = 6845274
$ws.Cells.Item(50, 5).Value = 'Partick'
$ws.Cells.Item(50, 6).Value = 'Ayr'
$ws.Cells.Item(50, 7).Value = 2
$ws.Cells.Item(50, 9).Value = 1
$ws.Cells.Item(50, 11).Value = 'D'
$ws.Cells.Item(50, 12).Value = 1.8
$ws.Cells.Item(50, 14).Value = 4
$ws.Cells.Item(50, 15).Value = 1.909
$ws.Cells.Item(50, 16).Value = 3.6
$ws.Cells.Item(50, 17).Value = 3.4
$ws.Cells.Item(50, 18).Value = -0.5
$ws.Cells.Item(50, 19).Value = 1.95
$ws.Cells.Item(50, 20).Value = 1.85
$ws.Cells.Item(50, 21).Value = 3
$ws.Cells.Item(50, 22).Value = 2
$ws.Cells.Item(50, 23).Value = 1.8
$ws.Cells.Item(50, 25).Value = 2.6
$ws.Cells.Item(50, 26).Value = -1
$ws.Cells.Item(50, 28).Value = 0.8500000000000001
$ws.Cells.Item(50, 29).Value = 1
$ws.Cells.Item(50, 30).Value = -1

# Row 51
$ws.Cells.Item(51, 2).Value = 6845270
$ws.Cells.Item(51, 5).Value = 'Airdrieonians'
$ws.Cells.Item(51, 6).Value = 'Dundee Utd'
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 'A'
$ws.Cells.Item(51, 12).Value = 3.75
$ws.Cells.Item(51, 14).Value = 1.85
$ws.Cells.Item(51, 15).Value = 5.25
$ws.Cells.Item(51, 16).Value = 4
$ws.Cells.Item(51, 17).Value = 1.533
$ws.Cells.Item(51, 18).Value = 1
$ws.Cells.Item(51, 19).Value = 1.9
$ws.Cells.Item(51, 20).Value = 1.95
$ws.Cells.Item(51, 21).Value = 2.75
$ws.Cells.Item(51, 22).Value = 1.85
$ws.Cells.Item(51, 23).Value = 2
$ws.Cells.Item(51, 25).Value = -1
$ws.Cells.Item(51, 26).Value = 0.5329999999999999
$ws.Cells.Item(51, 28).Value = 0.95
$ws.Cells.Item(51, 29).Value = -1
$ws.Cells.Item(51, 30).Value = 1

# Row 85
$ws.Cells.Item(85, 2).Value = 6845308
$ws.Cells.Item(85, 5).Value = 'Queens Park'
$ws.Cells.Item(85, 6).Value = 'Airdrieonians'
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 11).Value = 'A'
$ws.Cells.Item(85, 12).Value = 2.45
$ws.Cells.Item(85, 13).Value = 3.4
$ws.Cells.Item(85, 14).Value = 2.375
$ws.Cells.Item(85, 15).Value = 2.4
$ws.Cells.Item(85, 16).Value = 3.4
$ws.Cells.Item(85, 17).Value = 2.4
$ws.Cells.Item(85, 18).Value = 0
$ws.Cells.Item(85, 19).Value = 1.9
$ws.Cells.Item(85, 20).Value = 1.9
$ws.Cells.Item(85, 21).Value = 2.5
$ws.Cells.Item(85, 22).Value = 1.875
$ws.Cells.Item(85, 23).Value = 1.925
$ws.Cells.Item(85, 25).Value = -1
$ws.Cells.Item(85, 26).Value = 1.4
$ws.Cells.Item(85, 28).Value = 0.8999999999999999
$ws.Cells.Item(85, 29).Value = 0.875

# Row 86
$ws.Cells.Item(86, 2).Value = 6845309
$ws.Cells.Item(86, 5).Value = 'Raith'
$ws.Cells.Item(86, 6).Value = 'Arbroath'
$ws.Cells.Item(86, 7).Value = 2
$ws.Cells.Item(86, 8).Value = 2
$ws.Cells.Item(86, 9).Value = 1
$ws.Cells.Item(86, 12).Value = 1.444
$ws.Cells.Item(86, 13).Value = 4.2
$ws.Cells.Item(86, 14).Value = 5.5
$ws.Cells.Item(86, 15).Value = 1.4
$ws.Cells.Item(86, 16).Value = 4.333
$ws.Cells.Item(86, 17).Value = 6
$ws.Cells.Item(86, 18).Value = -1.25
$ws.Cells.Item(86, 19).Value = 1.875
$ws.Cells.Item(86, 20).Value = 1.925
$ws.Cells.Item(86, 21).Value = 3
$ws.Cells.Item(86, 22).Value = 1.9
$ws.Cells.Item(86, 25).Value = 3.333
$ws.Cells.Item(86, 27).Value = -1
$ws.Cells.Item(86, 28).Value = 0.925
$ws.Cells.Item(86, 29).Value = 0.8999999999999999
$ws.Cells.Item(86, 30).Value = -1

# Row 87
$ws.Cells.Item(87, 2).Value = 6845307
$ws.Cells.Item(87, 5).Value = 'Inverness CT'
$ws.Cells.Item(87, 6).Value = 'Morton'
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 11).Value = 'D'
$ws.Cells.Item(87, 12).Value = 2.2
$ws.Cells.Item(87, 14).Value = 2.75
$ws.Cells.Item(87, 15).Value = 2.1
$ws.Cells.Item(87, 16).Value = 3.1
$ws.Cells.Item(87, 17).Value = 3.1
$ws.Cells.Item(87, 18).Value = -0.25
$ws.Cells.Item(87, 19).Value = 2.025
$ws.Cells.Item(87, 20).Value = 1.825
$ws.Cells.Item(87, 21).Value = 2
$ws.Cells.Item(87, 22).Value = 1.95
$ws.Cells.Item(87, 23).Value = 1.9
$ws.Cells.Item(87, 25).Value = 2.1
$ws.Cells.Item(87, 26).Value = -1
$ws.Cells.Item(87, 27).Value = -0.5
$ws.Cells.Item(87, 28).Value = 0.4125
$ws.Cells.Item(87, 29).Value = -1
$ws.Cells.Item(87, 30).Value = 0.8999999999999999

# Row 95
$ws.Cells.Item(95, 2).Value = 6845316
$ws.Cells.Item(95, 5).Value = 'Ayr'
$ws.Cells.Item(95, 6).Value = 'Inverness CT'
$ws.Cells.Item(95, 7).Value = 1
$ws.Cells.Item(95, 9).Value = 0
$ws.Cells.Item(95, 10).Value = 1
$ws.Cells.Item(95, 12).Value = 2.25
$ws.Cells.Item(95, 13).Value = 3.75
$ws.Cells.Item(95, 14).Value = 2.5
$ws.Cells.Item(95, 15).Value = 2.15
$ws.Cells.Item(95, 16).Value = 3.75
$ws.Cells.Item(95, 17).Value = 2.8
$ws.Cells.Item(95, 18).Value = -0.25
$ws.Cells.Item(95, 19).Value = 1.975
$ws.Cells.Item(95, 20).Value = 1.825
$ws.Cells.Item(95, 21).Value = 2.25
$ws.Cells.Item(95, 22).Value = 1.85
$ws.Cells.Item(95, 23).Value = 1.95
$ws.Cells.Item(95, 26).Value = 1.8
$ws.Cells.Item(95, 28).Value = 0.825
$ws.Cells.Item(95, 29).Value = 0.8500000000000001

# Row 97
$ws.Cells.Item(97, 2).Value = 6845317
$ws.Cells.Item(97, 5).Value = 'Dundee Utd'
$ws.Cells.Item(97, 6).Value = 'Morton'
$ws.Cells.Item(97, 7).Value = 2
$ws.Cells.Item(97, 9).Value = 2
$ws.Cells.Item(97, 10).Value = 2
$ws.Cells.Item(97, 12).Value = 1.25
$ws.Cells.Item(97, 13).Value = 6
$ws.Cells.Item(97, 14).Value = 7
$ws.Cells.Item(97, 15).Value = 1.4
$ws.Cells.Item(97, 16).Value = 5.5
$ws.Cells.Item(97, 17).Value = 5.25
$ws.Cells.Item(97, 18).Value = -1.25
$ws.Cells.Item(97, 19).Value = 1.9
$ws.Cells.Item(97, 20).Value = 1.9
$ws.Cells.Item(97, 21).Value = 2.5
$ws.Cells.Item(97, 22).Value = 1.8
$ws.Cells.Item(97, 23).Value = 2
$ws.Cells.Item(97, 26).Value = 4.25
$ws.Cells.Item(97, 28).Value = 0.8999999999999999
$ws.Cells.Item(97, 29).Value = 0.8

# Row 160
$ws.Cells.Item(160, 2).Value = 7939128
$ws.Cells.Item(160, 5).Value = 'Ayr'
$ws.Cells.Item(160, 6).Value = 'Morton'
$ws.Cells.Item(160, 8).Value = 1
$ws.Cells.Item(160, 9).Value = 1
$ws.Cells.Item(160, 10).Value = 1
$ws.Cells.Item(160, 11).Value = 'D'
$ws.Cells.Item(160, 12).Value = 2.5
$ws.Cells.Item(160, 13).Value = 3.4
$ws.Cells.Item(160, 14).Value = 2.5
$ws.Cells.Item(160, 15).Value = 2.6
$ws.Cells.Item(160, 16).Value = 3.25
$ws.Cells.Item(160, 17).Value = 2.55
$ws.Cells.Item(160, 18).Value = 0
$ws.Cells.Item(160, 19).Value = 1.95
$ws.Cells.Item(160, 20).Value = 1.85
$ws.Cells.Item(160, 22).Value = 1.95
$ws.Cells.Item(160, 23).Value = 1.85
$ws.Cells.Item(160, 25).Value = 2.25
$ws.Cells.Item(160, 26).Value = -1
$ws.Cells.Item(160, 27).Value = 0
$ws.Cells.Item(160, 28).Value = 0
$ws.Cells.Item(160, 29).Value = -0.5
$ws.Cells.Item(160, 30).Value = 0.425

# Row 161
$ws.Cells.Item(161, 2).Value = 7842490
$ws.Cells.Item(161, 5).Value = 'Raith'
$ws.Cells.Item(161, 6).Value = 'Airdrieonians'
$ws.Cells.Item(161, 8).Value = 3
$ws.Cells.Item(161, 9).Value = 0
$ws.Cells.Item(161, 10).Value = 2
$ws.Cells.Item(161, 11).Value = 'A'
$ws.Cells.Item(161, 12).Value = 1.833
$ws.Cells.Item(161, 13).Value = 3.5
$ws.Cells.Item(161, 14).Value = 3.75
$ws.Cells.Item(161, 15).Value = 2
$ws.Cells.Item(161, 16).Value = 3.2
$ws.Cells.Item(161, 17).Value = 3.5
$ws.Cells.Item(161, 18).Value = -0.25
$ws.Cells.Item(161, 19).Value = 1.8
$ws.Cells.Item(161, 20).Value = 2
$ws.Cells.Item(161, 22).Value = 2
$ws.Cells.Item(161, 23).Value = 1.8
$ws.Cells.Item(161, 25).Value = -1
$ws.Cells.Item(161, 26).Value = 2.5
$ws.Cells.Item(161, 27).Value = -1
$ws.Cells.Item(161, 28).Value = 1
$ws.Cells.Item(161, 29).Value = 1
$ws.Cells.Item(161, 30).Value = -1

